# Apply the "break out stock.yaml completed" edit to the "week" sheet:
#  1. Convert the bsecode values in D57:D61 from text to real numbers
#     (same digits, now stored as numeric cells).
#  2. Append five new rows (62-66) duplicating the nsecode/name/bsecode of
#     rows 57-61 with a fresh per_chg/close/volume/timestamp snapshot,
#     keeping bsecode as text the way the earlier rows originally were.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# --- 1. D57:D61 inlineStr -> numeric -------------------------------------
$bseCodes = @{ 57 = 532830; 58 = 532296; 59 = 532400; 60 = 532482; 61 = 500049 }
foreach ($row in $bseCodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $bseCodes[$row]
}

# --- 2. New rows 62-66 ----------------------------------------------------
$newRows = @(
    @{ Row = 62; Sr = 1; NseCode = "ASTRAL";   Name = "Astral Poly Technik Limited";      Bse = "532830"; PerChg = 0.32;  Close = 2211.55; Volume = 314667;   Stamp = "20/06/2024 11:38:11" },
    @{ Row = 63; Sr = 2; NseCode = "GLENMARK"; Name = "Glenmark Pharmaceuticals Limited"; Bse = "532296"; PerChg = 0.45;  Close = 1241.25; Volume = 310556;   Stamp = "20/06/2024 11:38:11" },
    @{ Row = 64; Sr = 3; NseCode = "BSOFT";    Name = "Birlasoft Ltd";                    Bse = "532400"; PerChg = -0.42; Close = 687.95;  Volume = 3307590;  Stamp = "20/06/2024 11:38:11" },
    @{ Row = 65; Sr = 4; NseCode = "GRANULES"; Name = "Granules India Limited";           Bse = "532482"; PerChg = 0.73;  Close = 472.3;   Volume = 691935;   Stamp = "20/06/2024 11:38:11" },
    @{ Row = 66; Sr = 5; NseCode = "BEL";      Name = "Bharat Electronics Limited";       Bse = "500049"; PerChg = 0.84;  Close = 311.9;   Volume = 39505506; Stamp = "20/06/2024 11:38:11" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Sr
    $ws.Cells.Item($row, 2).Value = $r.NseCode
    $ws.Cells.Item($row, 3).Value = $r.Name

    # Keep bsecode stored as text (matches how rows 57-61 originally looked
    # before this edit converted them to numbers).
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $r.Bse
    $cell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $r.PerChg
    $ws.Cells.Item($row, 6).Value = $r.Close
    $ws.Cells.Item($row, 7).Value = $r.Volume
    $ws.Cells.Item($row, 8).Value = "week"
    $ws.Cells.Item($row, 9).Value = $r.Stamp
}
